# Insert a new data row at row 78 (pushes existing rows 78-196 down to 79-197)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(78).Insert()

$ws.Cells.Item(78, 1).Value = 10
$ws.Cells.Item(78, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(78, 3).Value = "La Araucanía"
$ws.Cells.Item(78, 4).Value = 44495
$ws.Cells.Item(78, 5).Value = 9
$ws.Cells.Item(78, 6).Value = 100112017
$ws.Cells.Item(78, 7).Value = "Apio"
$ws.Cells.Item(78, 8).Value = "Americana (o)"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 80
$ws.Cells.Item(78, 11).Value = 9000
$ws.Cells.Item(78, 12).Value = 9000
$ws.Cells.Item(78, 13).Value = 9000
$ws.Cells.Item(78, 14).Value = "$/docena de matas"
$ws.Cells.Item(78, 15).Value = "Región Metropolitana"
$ws.Cells.Item(78, 16).Value = 1500
$ws.Cells.Item(78, 17).Value = 6
$ws.Cells.Item(78, 18).Value = "Hortaliza"
